$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update cryptocurrency price (D) and 1h volume change (E) columns
# Cells whose new value looks like a plain number need a leading
# apostrophe so Excel keeps storing them as text (matching the
# original inlineStr cells), then the style is reset to Normal so
# no stray number-format/quote-prefix style is left on the cell.

$ws.Range('D2').Value = '67.745.78'
$ws.Range('E2').Value = '  -1.74%  '
$ws.Range('D3').Value = '2.674.68'
$ws.Range('E3').Value = '  -2.04%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').Value = '''599.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.31%  '
$ws.Range('D6').Value = '''167.19'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.55%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('D9').Value = '2.674.34'
$ws.Range('E9').Value = '  -2.08%  '
$ws.Range('E10').Value = '  +1.34%  '
$ws.Range('E11').Value = '  +1.01%  '
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('E13').Value = '  -1.63%  '
$ws.Range('E14').Value = '  -2.22%  '
$ws.Range('E15').Value = '  -1.89%  '
$ws.Range('E16').Value = '  -2.78%  '
$ws.Range('D17').Value = '67.691.55'
$ws.Range('E17').Value = '  -1.58%  '
$ws.Range('D18').Value = '2.676.44'
$ws.Range('E18').Value = '  -2.09%  '
$ws.Range('D19').Value = '''11.77'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.26%  '
$ws.Range('D20').Value = '''7.85'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.66%  '
$ws.Range('D21').Value = '''364.78'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.42%  '
$ws.Range('D22').Value = '''4.40'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.35%  '
$ws.Range('E23').Value = '  -2.20%  '
$ws.Range('E24').Value = '  -3.86%  '
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('D26').Value = '''70.84'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.16%  '
$ws.Range('D27').Value = '''10.21'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.80%  '
$ws.Range('E28').Value = '  -1.55%  '
$ws.Range('E29').Value = '  -3.21%  '
$ws.Range('D30').Value = '''0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.30%  '
$ws.Range('D31').Value = '''556.72'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.41%  '
$ws.Range('D32').Value = '''8.02'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.36%  '
$ws.Range('E33').Value = '  -4.14%  '
$ws.Range('E34').Value = '  -1.65%  '
$ws.Range('E35').Value = '  -1.54%  '
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('E37').Value = '  -4.68%  '
$ws.Range('D38').Value = '''19.50'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.59%  '
$ws.Range('D39').Value = '''154.67'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.67%  '
$ws.Range('D40').Value = '''0.373'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.45%  '
$ws.Range('E41').Value = '  -2.89%  '
$ws.Range('E42').Value = '  -4.40%  '
$ws.Range('D43').Value = '''17.95'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.31%  '
$ws.Range('D44').Value = '''2.53'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.80%  '
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').Value = '''40.36'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.94%  '
$ws.Range('D47').Value = '0.0₆0301'
$ws.Range('E47').Value = '  -4.96%  '
$ws.Range('E48').Value = '  -3.11%  '
$ws.Range('D49').Value = '''153.90'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.92%  '
$ws.Range('D50').Value = '''3.88'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.98%  '
$ws.Range('E51').Value = '  -3.68%  '
